# Auto-generated script: apply 2024-08-31 data update to violent-crime-full-year.xlsx
# Updates year-2024 (column K, occasionally H/J corrections) values across 43 worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5328
$ws.Range("J3").Value = 8077
$ws.Range("K3").Value = 5495
$ws.Range("H4").Value = 1741
$ws.Range("K4").Value = 1143
$ws.Range("K5").Value = 392
$ws.Range("K6").Value = 6101
$ws.Range("H7").Value = 26054
$ws.Range("J7").Value = 29297
$ws.Range("K7").Value = 18459

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 339
$ws.Range("K7").Value = 1235

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 138
$ws.Range("K3").Value = 149
$ws.Range("K6").Value = 93
$ws.Range("K7").Value = 410

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 216
$ws.Range("K3").Value = 292
$ws.Range("H4").Value = 62
$ws.Range("K6").Value = 232
$ws.Range("H7").Value = 1303
$ws.Range("K7").Value = 792

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 314

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K3").Value = 106
$ws.Range("K7").Value = 418

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 131
$ws.Range("K7").Value = 313

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 162
$ws.Range("K7").Value = 548
$ws.Range("K8").Value = 1235
$ws.Range("K11").Value = 351
$ws.Range("K14").Value = 97
$ws.Range("K19").Value = 543
$ws.Range("K20").Value = 429
$ws.Range("K21").Value = 58
$ws.Range("K22").Value = 49
$ws.Range("K26").Value = 24
$ws.Range("K29").Value = 984
$ws.Range("K30").Value = 72
$ws.Range("J31").Value = 304
$ws.Range("K31").Value = 202
$ws.Range("H33").Value = 1303
$ws.Range("K33").Value = 792
$ws.Range("K36").Value = 246
$ws.Range("K42").Value = 684
$ws.Range("K43").Value = 164
$ws.Range("K46").Value = 38
$ws.Range("K47").Value = 125
$ws.Range("K49").Value = 103
$ws.Range("K51").Value = 233
$ws.Range("K52").Value = 482
$ws.Range("K54").Value = 361
$ws.Range("K55").Value = 205
$ws.Range("K57").Value = 69
$ws.Range("K63").Value = 52
$ws.Range("K64").Value = 118
$ws.Range("K65").Value = 418
$ws.Range("K67").Value = 700
$ws.Range("K73").Value = 159
$ws.Range("K76").Value = 255
$ws.Range("K78").Value = 212
$ws.Range("K79").Value = 458
$ws.Range("K83").Value = 410
$ws.Range("K84").Value = 138
$ws.Range("K85").Value = 872
$ws.Range("K89").Value = 270
$ws.Range("K90").Value = 166
$ws.Range("K91").Value = 203
$ws.Range("K95").Value = 314
$ws.Range("K96").Value = 200
$ws.Range("K99").Value = 313
$ws.Range("H101").Value = 26054
$ws.Range("J101").Value = 29297
$ws.Range("K101").Value = 18459

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 68
$ws.Range("J3").Value = 71
$ws.Range("J7").Value = 304
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 201
$ws.Range("K3").Value = 250
$ws.Range("K6").Value = 195
$ws.Range("K7").Value = 700

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K2").Value = 44
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 361

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 282
$ws.Range("K3").Value = 356
$ws.Range("K7").Value = 984

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 159
$ws.Range("K3").Value = 172
$ws.Range("K7").Value = 543

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 34
$ws.Range("K7").Value = 97

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 180
$ws.Range("K3").Value = 212
$ws.Range("K6").Value = 259
$ws.Range("K7").Value = 684

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 200

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 51
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K4").Value = 31
$ws.Range("K7").Value = 458

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 142
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 429

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 96
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 246

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 186
$ws.Range("K4").Value = 19
$ws.Range("K6").Value = 146
$ws.Range("K7").Value = 548

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K6").Value = 42
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 351

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 76
$ws.Range("K4").Value = 32
$ws.Range("K7").Value = 270

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 52
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 67
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 164

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 291
$ws.Range("K3").Value = 294
$ws.Range("K6").Value = 211
$ws.Range("K7").Value = 872

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 128
$ws.Range("K6").Value = 176
$ws.Range("K7").Value = 482

Write-Output "Updated 165 cells across 43 worksheets"